$d = $word.ActiveDocument

# The three <id>...</id> tags were each split across 3 separate runs
# (open-tag run, id-value run, close-tag run). Collapse each trio into a
# single run containing the full "<id>pNNNv_N</id>" text, keeping the
# Courier New / gold-colored formatting of the tag runs.
$ids = @("p134v_1", "p134v_2", "p134v_3")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $new = "<id>" + $id + "</id>"
    $found = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Host "Replaced $id : $found"
}
